$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update fresnelEffectiveDiameter value (B15): 100 -> 125
$ws.Range("B15").Value = 125

# Update focalLength value (B16): 3 -> 6  (longer length, used in fresnel)
$ws.Range("B16").Value = 6

# Move active selection to B16
$ws.Range("B16").Select()
